$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skills")

# Clear out the old data rows (rows 2-35); header row (row 1) stays the same.
$ws.Range("A2:D35").Clear()

# Build the new skills table (32 data rows x 4 columns) to replace rows 2-33.
$data = New-Object 'object[,]' 32,4
$data[0,0] = 'Tech Stack'
$data[0,1] = 'Coding Languages'
$data[0,2] = 'R Statistical programming language'
$data[0,3] = '5 years'
$data[1,0] = 'Tech Stack'
$data[1,1] = 'Other'
$data[1,2] = 'Markdown'
$data[1,3] = '5 years'
$data[2,0] = 'Tech Stack'
$data[2,1] = 'Other'
$data[2,2] = 'Git'
$data[2,3] = '5 years'
$data[3,0] = 'Tech Stack'
$data[3,1] = 'Software'
$data[3,2] = 'QGIS'
$data[3,3] = '7 years'
$data[4,0] = 'Tech Stack'
$data[4,1] = 'Software'
$data[4,2] = 'ArcGIS'
$data[4,3] = '3 years'
$data[5,0] = 'Tech Stack'
$data[5,1] = 'Coding Languages'
$data[5,2] = 'SQL'
$data[5,3] = ''
$data[6,0] = 'Tech Stack'
$data[6,1] = 'Coding Languages'
$data[6,2] = 'Python'
$data[6,3] = ''
$data[7,0] = 'Tech Stack'
$data[7,1] = 'Coding Languages'
$data[7,2] = 'HTML'
$data[7,3] = ''
$data[8,0] = 'Tech Stack'
$data[8,1] = 'Coding Languages'
$data[8,2] = 'CSS'
$data[8,3] = ''
$data[9,0] = 'Tech Stack'
$data[9,1] = 'Software'
$data[9,2] = 'Tableau'
$data[9,3] = ''
$data[10,0] = 'Language'
$data[10,1] = ''
$data[10,2] = 'Hebrew'
$data[10,3] = 'C2 - Native'
$data[11,0] = 'Language'
$data[11,1] = ''
$data[11,2] = 'English'
$data[11,3] = 'C2 - Fluent'
$data[12,0] = 'Language'
$data[12,1] = ''
$data[12,2] = 'Spanish'
$data[12,3] = 'B1 - Intermediate'
$data[13,0] = 'Tech Stack'
$data[13,1] = 'Software'
$data[13,2] = 'Google Sheets'
$data[13,3] = ''
$data[14,0] = 'Tech Stack'
$data[14,1] = 'Coding Languages'
$data[14,2] = 'Shiny'
$data[14,3] = ''
$data[15,0] = 'Tech Stack'
$data[15,1] = 'Software'
$data[15,2] = 'Canva'
$data[15,3] = ''
$data[16,0] = 'Tech Stack'
$data[16,1] = 'Software'
$data[16,2] = 'Adobe Photoshop'
$data[16,3] = ''
$data[17,0] = 'Tech Stack'
$data[17,1] = 'Software'
$data[17,2] = 'Adobe Illustrator'
$data[17,3] = ''
$data[18,0] = 'Tech Stack'
$data[18,1] = 'Software'
$data[18,2] = 'GIMP'
$data[18,3] = ''
$data[19,0] = 'Tech Stack'
$data[19,1] = 'Software'
$data[19,2] = 'InkScape'
$data[19,3] = ''
$data[20,0] = 'Tech Stack'
$data[20,1] = 'Software'
$data[20,2] = 'Microsoft Office Suite'
$data[20,3] = ''
$data[21,0] = 'Tech Stack'
$data[21,1] = 'Software'
$data[21,2] = 'Airtable'
$data[21,3] = ''
$data[22,0] = 'Tech Stack'
$data[22,1] = 'Software'
$data[22,2] = 'Notion'
$data[22,3] = ''
$data[23,0] = 'Tech Stack'
$data[23,1] = 'Software'
$data[23,2] = 'Trello'
$data[23,3] = ''
$data[24,0] = 'Soft Skills'
$data[24,1] = ''
$data[24,2] = 'Strong communication skills'
$data[24,3] = ''
$data[25,0] = 'Soft Skills'
$data[25,1] = ''
$data[25,2] = 'Experience with data management tools and processes'
$data[25,3] = ''
$data[26,0] = 'Soft Skills'
$data[26,1] = ''
$data[26,2] = 'Familiarity with data piplines'
$data[26,3] = ''
$data[27,0] = 'Soft Skills'
$data[27,1] = ''
$data[27,2] = 'Demonstrated analytical and problem-solving skills'
$data[27,3] = ''
$data[28,0] = 'Soft Skills'
$data[28,1] = ''
$data[28,2] = 'Statistical knowledge'
$data[28,3] = ''
$data[29,0] = 'Soft Skills'
$data[29,1] = ''
$data[29,2] = 'Ability to work independently'
$data[29,3] = ''
$data[30,0] = 'Soft Skills'
$data[30,1] = ''
$data[30,2] = 'Public speaking experience'
$data[30,3] = ''
$data[31,0] = 'Soft Skills'
$data[31,1] = ''
$data[31,2] = 'Data storytelling'
$data[31,3] = ''

$ws.Range("A2:D33").Value = $data

# Match the author's final cursor position recorded in the saved workbook.
$ws.Range("B27").Select() | Out-Null
